$wb = $excel.ActiveWorkbook

# Sheets that contain the "想去人数" (wish-to-attend count) column needing updates:
# "展览" (sheet1) and "全部类型" (sheet4) contain identical data rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 343
    $ws.Range("F4").Value = 1519
    $ws.Range("F9").Value = 336
}
